# Updated cryptos list on Fri Jul  5 04:14:02 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# re-sorts the Maker/Stacks pair (rows 47-48) to reflect the new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then new values for the B/C/D/E columns that changed.
# `$null` means that column is left untouched for that row.
$updates = @(
    @{ Row = 2; B = $null; C = $null; D = "55.270.45"; E = "  -6.40%  " }
    @{ Row = 3; B = $null; C = $null; D = "2.918.27"; E = "  -10.00%  " }
    @{ Row = 5; B = $null; C = $null; D = "466.76"; E = "  -13.67%  " }
    @{ Row = 6; B = $null; C = $null; D = "123.60"; E = "  -9.66%  " }
    @{ Row = 7; B = $null; C = $null; D = $null; E = "  -0.03%  " }
    @{ Row = 8; B = $null; C = $null; D = "2.918.24"; E = "  -9.98%  " }
    @{ Row = 9; B = $null; C = $null; D = "0.398"; E = "  -13.35%  " }
    @{ Row = 10; B = $null; C = $null; D = "6.49"; E = "  -15.09%  " }
    @{ Row = 11; B = $null; C = $null; D = "0.0945"; E = "  -18.03%  " }
    @{ Row = 12; B = $null; C = $null; D = "0.323"; E = "  -18.56%  " }
    @{ Row = 13; B = $null; C = $null; D = $null; E = "  -4.18%  " }
    @{ Row = 14; B = $null; C = $null; D = "3.420.86"; E = "  -9.94%  " }
    @{ Row = 15; B = $null; C = $null; D = "22.34"; E = "  -14.27%  " }
    @{ Row = 16; B = $null; C = $null; D = "55.255.85"; E = "  -6.52%  " }
    @{ Row = 17; B = $null; C = $null; D = "2.923.53"; E = "  -9.81%  " }
    @{ Row = 18; B = $null; C = $null; D = $null; E = "  -17.62%  " }
    @{ Row = 19; B = $null; C = $null; D = "5.02"; E = "  -15.05%  " }
    @{ Row = 20; B = $null; C = $null; D = "11.30"; E = "  -14.68%  " }
    @{ Row = 21; B = $null; C = $null; D = "6.94"; E = "  -16.25%  " }
    @{ Row = 22; B = $null; C = $null; D = "304.30"; E = "  -16.05%  " }
    @{ Row = 23; B = $null; C = $null; D = "1.00"; E = "  -0.07%  " }
    @{ Row = 24; B = $null; C = $null; D = "0.438"; E = "  -15.89%  " }
    @{ Row = 25; B = $null; C = $null; D = "58.61"; E = "  -16.83%  " }
    @{ Row = 26; B = $null; C = $null; D = "0.998"; E = "  -0.24%  " }
    @{ Row = 27; B = $null; C = $null; D = "0.151"; E = "  -11.35%  " }
    @{ Row = 28; B = $null; C = $null; D = $null; E = "  -0.12%  " }
    @{ Row = 29; B = $null; C = $null; D = $null; E = "  -19.15%  " }
    @{ Row = 30; B = $null; C = $null; D = "5.86"; E = "  -16.95%  " }
    @{ Row = 31; B = $null; C = $null; D = $null; E = "  -11.06%  " }
    @{ Row = 32; B = $null; C = $null; D = $null; E = "  -14.57%  " }
    @{ Row = 33; B = $null; C = $null; D = $null; E = "  -16.08%  " }
    @{ Row = 34; B = $null; C = $null; D = "1.56"; E = "  -19.65%  " }
    @{ Row = 35; B = $null; C = $null; D = "141.25"; E = "  -13.85%  " }
    @{ Row = 36; B = $null; C = $null; D = "4.10"; E = "  -17.03%  " }
    @{ Row = 37; B = $null; C = $null; D = "5.29"; E = "  -17.57%  " }
    @{ Row = 38; B = $null; C = $null; D = "1.18"; E = "  -17.76%  " }
    @{ Row = 39; B = $null; C = $null; D = "2.948.02"; E = "  -9.98%  " }
    @{ Row = 40; B = $null; C = $null; D = $null; E = "  -0.02%  " }
    @{ Row = 41; B = $null; C = $null; D = "0.0600"; E = "  -15.57%  " }
    @{ Row = 42; B = $null; C = $null; D = "21.16"; E = "  -19.41%  " }
    @{ Row = 43; B = $null; C = $null; D = "34.78"; E = "  -15.67%  " }
    @{ Row = 44; B = $null; C = $null; D = "0.942"; E = "  -14.66%  " }
    @{ Row = 45; B = $null; C = $null; D = "0.598"; E = "  -16.94%  " }
    @{ Row = 46; B = $null; C = $null; D = "3.34"; E = "  -17.14%  " }
    @{ Row = 47; B = "Maker"; C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D = "2.023.80"; E = "  -11.99%  " }
    @{ Row = 48; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "1.28"; E = "  -15.29%  " }
    @{ Row = 49; B = $null; C = $null; D = "5.16"; E = "  -18.31%  " }
    @{ Row = 50; B = $null; C = $null; D = "17.25"; E = "  -17.41%  " }
    @{ Row = 51; B = $null; C = $null; D = "0.0206"; E = "  -15.15%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.B) { $ws.Cells.Item($row, 2).Value = $u.B }
    if ($null -ne $u.C) { $ws.Cells.Item($row, 3).Value = $u.C }
    if ($null -ne $u.D) {
        # Price column holds text that often *looks* numeric (e.g. "1.00",
        # "0.0600", "5.86"). Force text so Excel does not silently convert it
        # to a number, then clear the temporary format so the cell keeps the
        # workbook default style (matches the source data, which is plain text).
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    if ($null -ne $u.E) { $ws.Cells.Item($row, 5).Value = $u.E }
}
